$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts existing D:K data to F:M)
$ws.Columns("D:E").Insert()

# Copy number formatting from the (now-shifted) F:G columns into new D:E so styles match
$ws.Columns("F:G").Copy()
$ws.Columns("D:E").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarter columns (D, E) with the new data for each row
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 318100
$ws.Range("E8").Value = 406600
$ws.Range("D9").Value = 245700
$ws.Range("E9").Value = 311200
$ws.Range("D10").Value = 72400
$ws.Range("E10").Value = 95400
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 2000
$ws.Range("E15").Value = 2000
$ws.Range("D17").Value = 293200
$ws.Range("E17").Value = 359800
$ws.Range("D18").Value = 24900
$ws.Range("E18").Value = 46800
$ws.Range("D20").Value = -600
$ws.Range("E20").Value = -100
$ws.Range("D21").Value = 41800
$ws.Range("E21").Value = 64200
$ws.Range("D22").Value = 5700
$ws.Range("E22").Value = 4500
$ws.Range("D23").Value = 18600
$ws.Range("E23").Value = 42100
$ws.Range("D24").Value = -2100
$ws.Range("E24").Value = 12200
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 20700
$ws.Range("E26").Value = 29900
$ws.Range("D27").Value = 18900
$ws.Range("E27").Value = 25900
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -4600
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 600
$ws.Range("E32").Value = 100
$ws.Range("D33").Value = 14300
$ws.Range("E33").Value = 25900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 14300
$ws.Range("E35").Value = 25900
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 19800
$ws.Range("E41").Value = 17600
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 152500
$ws.Range("E43").Value = 237000
$ws.Range("D44").Value = 246500
$ws.Range("E44").Value = 247200
$ws.Range("D45").Value = 7600
$ws.Range("E45").Value = 7200
$ws.Range("D46").Value = 426400
$ws.Range("E46").Value = 509000
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 402800
$ws.Range("E48").Value = 402200
$ws.Range("D49").Value = 141000
$ws.Range("E49").Value = 143700
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 36900
$ws.Range("E52").Value = 39300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1007000
$ws.Range("E54").Value = 1094200
$ws.Range("D57").Value = 66400
$ws.Range("E57").Value = 95900
$ws.Range("D58").Value = 49500
$ws.Range("E58").Value = 49400
$ws.Range("D59").Value = 74200
$ws.Range("E59").Value = 75600
$ws.Range("D60").Value = 190100
$ws.Range("E60").Value = 220800
$ws.Range("D61").Value = 264300
$ws.Range("E61").Value = 314800
$ws.Range("D62").Value = 57700
$ws.Range("E62").Value = 60500
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 524900
$ws.Range("E66").Value = 621800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 285100
$ws.Range("E70").Value = 285100
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 22000
$ws.Range("E72").Value = 11600
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 197000
$ws.Range("E76").Value = 187300
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 14300
$ws.Range("E81").Value = 25900
$ws.Range("D83").Value = 17500
$ws.Range("E83").Value = 17500
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 89900
$ws.Range("E89").Value = 48300
$ws.Range("D91").Value = -11800
$ws.Range("E91").Value = -12400
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -11200
$ws.Range("E94").Value = -11900
$ws.Range("D96").Value = -9500
$ws.Range("E96").Value = -5500
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -76300
$ws.Range("E100").Value = -37500
$ws.Range("D101").Value = -300
$ws.Range("E101").Value = 300
$ws.Range("D102").Value = 2100
$ws.Range("E102").Value = -800
